# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells stay as text so formats like "1.00" or "0.540"
# are not auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.227.33"
$ws.Range("E2").Value = "  +1.96%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.938.27"
$ws.Range("E3").Value = "  +2.18%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.30"
$ws.Range("E5").Value = "  +1.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.73"
$ws.Range("E6").Value = "  +5.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.936.56"
$ws.Range("E7").Value = "  +2.23%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("E9").Value = "  +1.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +2.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.46"
$ws.Range("E11").Value = "  +2.27%  "

$ws.Range("E12").Value = "  +3.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("E13").Value = "  +6.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.72"
$ws.Range("E14").Value = "  +5.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.601.31"
$ws.Range("E15").Value = "  +2.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.968.52"
$ws.Range("E16").Value = "  +2.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.225.79"
$ws.Range("E17").Value = "  +1.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.70"
$ws.Range("E18").Value = "  +2.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.63"
$ws.Range("E19").Value = "  +8.74%  "

$ws.Range("E20").Value = "  -0.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.08"
$ws.Range("E21").Value = "  -3.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "497.47"
$ws.Range("E22").Value = "  +2.77%  "

$ws.Range("E23").Value = "  +4.44%  "

$ws.Range("E24").Value = "  +5.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.13"
$ws.Range("E25").Value = "  +2.91%  "

$ws.Range("E26").Value = "  +2.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.42"
$ws.Range("E27").Value = "  +2.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.22"
$ws.Range("E28").Value = "  +2.17%  "

$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.02"
$ws.Range("E30").Value = "  +1.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.47"
$ws.Range("E31").Value = "  +3.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.090.51"
$ws.Range("E32").Value = "  +2.12%  "

$ws.Range("E33").Value = "  -0.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.43"
$ws.Range("E34").Value = "  +0.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.903.86"
$ws.Range("E35").Value = "  +2.66%  "

$ws.Range("E36").Value = "  +1.79%  "

$ws.Range("E37").Value = "  +5.31%  "

$ws.Range("E38").Value = "  +1.85%  "

$ws.Range("E39").Value = "  +1.62%  "

$ws.Range("E40").Value = "  +11.38%  "

$ws.Range("E41").Value = "  +4.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("E43").Value = "  +8.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "439.69"
$ws.Range("E44").Value = "  -0.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.37"
$ws.Range("E45").Value = "  -0.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.71"
$ws.Range("E46").Value = "  +3.88%  "

$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000278"
$ws.Range("E48").Value = "  +23.54%  "

$ws.Range("E49").Value = "  +3.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.82"
$ws.Range("E50").Value = "  +5.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.46"
$ws.Range("E51").Value = "  +0.51%  "

